$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 54
$ws.Range("Q54").Value = 1.75
$ws.Range("R54").Value = 2.05
$ws.Range("AD54").Value = 8.5
$ws.Range("AO54").Value = 7.5

# Row 56
$ws.Range("K56").Value = 2.3
$ws.Range("Q56").Value = 1.73
$ws.Range("R56").Value = 2.08
$ws.Range("S56").Value = 1.33
$ws.Range("T56").Value = 3.25
$ws.Range("AT56").Value = 3.25

# Row 61
$ws.Range("G61").Value = 2.18
$ws.Range("H61").Value = 2.95
$ws.Range("I61").Value = 3.4
$ws.Range("J61").Value = 2.82
$ws.Range("K61").Value = 1.95
$ws.Range("L61").Value = 4
$ws.Range("N61").Value = 6.45
$ws.Range("O61").Value = 1.42
$ws.Range("P61").Value = 2.47
$ws.Range("Q61").Value = 2.2
$ws.Range("S61").Value = 1.47
$ws.Range("T61").Value = 2.32
$ws.Range("U61").Value = 1.9
$ws.Range("V61").Value = 1.72
$ws.Range("W61").Value = 6.2
$ws.Range("X61").Value = 9.75
$ws.Range("Y61").Value = 9
$ws.Range("Z61").Value = 21
$ws.Range("AA61").Value = 20
$ws.Range("AB61").Value = 35
$ws.Range("AC61").Value = 7.1
$ws.Range("AD61").Value = 5.8
$ws.Range("AE61").Value = 16
$ws.Range("AF61").Value = 90
$ws.Range("AH61").Value = 8.25
$ws.Range("AI61").Value = 17
$ws.Range("AJ61").Value = 12
$ws.Range("AK61").Value = 50
$ws.Range("AL61").Value = 35
$ws.Range("AN61").Value = 3.9
$ws.Range("AO61").Value = 11.5
$ws.Range("AP61").Value = 22
$ws.Range("AQ61").Value = 50
$ws.Range("AT61").Value = 2.3
$ws.Range("AU61").Value = 7.3
$ws.Range("AW61").Value = 5.1
$ws.Range("AX61").Value = 20
$ws.Range("AZ61").Value = 110
$ws.Range("BA61").Value = 150

# Row 78
$ws.Range("G78").Value = 1.85
$ws.Range("I78").Value = 4
$ws.Range("U78").Value = 1.91
$ws.Range("V78").Value = 1.8
$ws.Range("W78").Value = 6.5
$ws.Range("X78").Value = 8.5
$ws.Range("AB78").Value = 29
$ws.Range("AF78").Value = 51
$ws.Range("AG78").Value = 900
$ws.Range("AL78").Value = 34
$ws.Range("AO78").Value = 11
$ws.Range("AQ78").Value = 41
$ws.Range("AU78").Value = 8.5
$ws.Range("BA78").Value = 101
